# Applies the "Updated cryptos list" data refresh to the crypto table.
# Each row 2..51 holds one coin: A=rank(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h).
# Values are written as literal text (matching the source workbook, which stores
# every data cell as text) rather than as numbers. For D-column values that Excel's
# own type-sniffing would otherwise parse as a number (losing the original text
# formatting, e.g. trailing zeros), a leading apostrophe forces a literal/text entry
# -- the same trick used when typing such values directly into Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.199.14'
$ws.Range('E2').Value = '  +1.21%  '

# Row 3
$ws.Range('D3').Value = '1.834.48'
$ws.Range('E3').Value = '  +1.05%  '

# Row 4
$ws.Range('E4').Value = '  +1.06%  '

# Row 5
$ws.Range('D5').Value = '''313.77'
$ws.Range('E5').Value = '  +1.26%  '

# Row 6
$ws.Range('D6').Value = '''1.010'
$ws.Range('E6').Value = '  +0.96%  '

# Row 7
$ws.Range('D7').Value = '''0.4713'
$ws.Range('E7').Value = '  +0.76%  '

# Row 8
$ws.Range('E8').Value = '  -0.14%  '

# Row 9
$ws.Range('D9').Value = '''0.07421'
$ws.Range('E9').Value = '  +0.71%  '

# Row 10
$ws.Range('D10').Value = '''0.8828'
$ws.Range('E10').Value = '  +1.56%  '

# Row 11
$ws.Range('E11').Value = '  +0.07%  '

# Row 12
$ws.Range('D12').Value = '1.826.92'
$ws.Range('E12').Value = '  +0.86%  '

# Row 13
$ws.Range('D13').Value = '''0.07324'
$ws.Range('E13').Value = '  +3.41%  '

# Row 14
$ws.Range('D14').Value = '''5.482'
$ws.Range('E14').Value = '  +2.24%  '

# Row 15
$ws.Range('D15').Value = '''92.84'
$ws.Range('E15').Value = '  +0.72%  '

# Row 16
$ws.Range('D16').Value = '''6.564'
$ws.Range('E16').Value = '  +1.00%  '

# Row 17
$ws.Range('D17').Value = '''1.013'
$ws.Range('E17').Value = '  +1.04%  '

# Row 18
$ws.Range('D18').Value = '''0.000008794'
$ws.Range('E18').Value = '  +0.82%  '

# Row 19
$ws.Range('E19').Value = '  +0.88%  '

# Row 20
$ws.Range('E20').Value = '  +0.51%  '

# Row 21
$ws.Range('D21').Value = '27.206.20'
$ws.Range('E21').Value = '  +1.09%  '

# Row 22
$ws.Range('D22').Value = '''5.304'
$ws.Range('E22').Value = '  -0.64%  '

# Row 23
$ws.Range('E23').Value = '  +1.36%  '

# Row 24
$ws.Range('D24').Value = '2.045.15'
$ws.Range('E24').Value = '  +0.31%  '

# Row 25
$ws.Range('E25').Value = '  +0.58%  '

# Row 26
$ws.Range('D26').Value = '''152.49'
$ws.Range('E26').Value = '  +0.51%  '

# Row 27
$ws.Range('E27').Value = '  +1.10%  '

# Row 28
$ws.Range('D28').Value = '''2.165'
$ws.Range('E28').Value = '  -0.88%  '

# Row 29
$ws.Range('D29').Value = '''5.277'
$ws.Range('E29').Value = '  -0.36%  '

# Row 30
$ws.Range('D30').Value = '''117.64'
$ws.Range('E30').Value = '  +1.86%  '

# Row 31
$ws.Range('D31').Value = '''0.08932'
$ws.Range('E31').Value = '  +0.10%  '

# Row 32
$ws.Range('D32').Value = '''0.7603'
$ws.Range('E32').Value = '  -0.73%  '

# Row 33
$ws.Range('E33').Value = '  +1.15%  '

# Row 34
$ws.Range('E34').Value = '  +1.61%  '

# Row 35
$ws.Range('D35').Value = '''2.943'
$ws.Range('E35').Value = '  +0.70%  '

# Row 36
$ws.Range('E36').Value = '  +0.97%  '

# Row 37
$ws.Range('E37').Value = '  +0.35%  '

# Row 38
$ws.Range('D38').Value = '''0.05338'
$ws.Range('E38').Value = '  +1.48%  '

# Row 39
$ws.Range('E39').Value = '  +0.27%  '

# Row 40
$ws.Range('D40').Value = '''3.010'
$ws.Range('E40').Value = '  +2.27%  '

# Row 41
$ws.Range('E41').Value = '  +2.90%  '

# Row 42
$ws.Range('D42').Value = '''7.346'
$ws.Range('E42').Value = '  +1.11%  '

# Row 43
$ws.Range('D43').Value = '''0.5355'
$ws.Range('E43').Value = '  +0.44%  '

# Row 44
$ws.Range('D44').Value = '''0.1664'
$ws.Range('E44').Value = '  +0.18%  '

# Row 45
$ws.Range('D45').Value = '''8.540'
$ws.Range('E45').Value = '  +1.23%  '

# Row 46
$ws.Range('D46').Value = '''0.4955'
$ws.Range('E46').Value = '  +0.15%  '

# Row 47
$ws.Range('E47').Value = '  +0.90%  '

# Row 48
$ws.Range('E48').Value = '  +1.02%  '

# Row 49
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''1.672'
$ws.Range('E49').Value = '  -0.01%  '

# Row 50
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''103.90'
$ws.Range('E50').Value = '  +1.14%  '

# Row 51
$ws.Range('E51').Value = '  +0.72%  '
